# vault backup: 2024-12-29 22:26:05
#
# The "Revenue Projection" sheet was a duplicate of the "Startup Budget"
# sheet's data. This backup clears that duplicated data back out (keeping
# the header/total row formatting in place) and leaves the user with the
# "Revenue Projection" tab selected/active, with cell G14 as the last
# selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Startup Budget")
$ws2 = $wb.Worksheets.Item("Revenue Projection")

# Wipe the budget-table rows that had been copied onto "Revenue Projection"
# (A2:C8) completely -- content AND formatting -- so those rows disappear
# from the sheet entirely, just like they never existed.
$ws2.Range("A2:C8").Clear()

# The header row and the two total rows keep their cell formatting
# (borders/bold/number format/alignment), only their values/formulas are
# removed.
$ws2.Range("A1:D1").ClearContents()
$ws2.Range("A9:C9").ClearContents()
$ws2.Range("A11:C11").ClearContents()

# Leave the user on the "Revenue Projection" tab with G14 as the last
# selected cell (so "Startup Budget" is no longer the active/selected tab).
$ws2.Activate()
$ws2.Range("G14").Select() | Out-Null
